# Apply the "Updated data from WRI for HK EPS 2.0 update" edit:
#   - BGDPbES sheet, row 5 ("hydro"): the 2015 Guaranteed Dispatch value (B5)
#     changes from 1 to 0. Every later year (C5:AK5) is driven by a shared
#     formula "=$B5", so they recompute to 0 automatically.
#   - The workbook was re-saved with the BGDPbES sheet as the active/selected
#     sheet (instead of the About sheet), with a new selection on each sheet.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("BGDPbES")

# --- data edit -------------------------------------------------------
# Hydro's 2015 guaranteed-dispatch percentage goes from 1 (100%) to 0.
$wsData.Range("B5").Value = 0

# --- view / selection state -------------------------------------------
# Leave a selection parked on the About sheet (matches the saved file, where
# "About" is no longer the tab shown when the workbook is reopened).
[void]$wsAbout.Range("A13").Select()

# Make BGDPbES the active/visible sheet, with its own selection restored.
[void]$wsData.Activate()
[void]$wsData.Range("AL17").Select()
